$d = $word.ActiveDocument

# Fix the typo "Im am studying EEE" -> "I am studying EEE"
$d.Content.Find.Execute("Im am studying EEE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "I am studying EEE", 2)

# Append a new paragraph with "I am from Tirupur" after the last paragraph,
# matching the formatting of the preceding text (Times New Roman, bold, 12pt, en-US).
$last = $d.Paragraphs.Last.Range
$last.InsertParagraphAfter()

$newRange = $d.Paragraphs.Last.Range
$newRange.Font.Name = "Times New Roman"
$newRange.Font.NameAscii = "Times New Roman"
$newRange.Font.Bold = $true
$newRange.Font.Size = 12
$newRange.InsertAfter("I am from Tirupur")
